$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planning & Journal")

# Row 51: Création des endpoints /projects - planned/realised hours updated
$ws.Range("D51").Value = 2
$ws.Range("E51").Value = 3

# Row 52: Création des endpoints /teams - planned/realised hours updated
$ws.Range("D52").Value = 4
$ws.Range("E52").Value = 4

# Row 53: realised hours updated
$ws.Range("E53").Value = 0

# Row 54: status switched from "A faire" to "Termine", realised hours updated
$ws.Range("C54").Value = "Terminé"
$ws.Range("E54").Value = 5

# Update selection to reflect where the user ended up after the edits
$ws.Range("G56").Select()
